# Applies the weekly fruit/vegetable price-data refresh for the
# "Hortaliza, Vega Monumental Concepcion - Arveja Verde" sheet:
# rows 2-9 are re-shuffled to their updated reporting-date order,
# carrying each row's Fecha/Variedad/Volumen/Precio*/Origen data along.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44483
$ws.Range("H2").Value = "Perfection"
$ws.Range("J2").Value = 220
$ws.Range("K2").Value = 19000
$ws.Range("L2").Value = 20000
$ws.Range("M2").Value = 19455
$ws.Range("N2").Value = "$/saco 25 kilos"
$ws.Range("O2").Value = "Región Metropolitana"
$ws.Range("P2").Value = 778

$ws.Range("D3").Value = 44162
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 17000
$ws.Range("L3").Value = 18000
$ws.Range("M3").Value = 17500
$ws.Range("N3").Value = "$/saco 25 kilos"
$ws.Range("O3").Value = "Región del Maule"
$ws.Range("P3").Value = 700

$ws.Range("D4").Value = 44335
$ws.Range("H4").Value = "Perfection"
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 30000
$ws.Range("L4").Value = 32000
$ws.Range("M4").Value = 31000
$ws.Range("N4").Value = "$/malla 25 kilos"
$ws.Range("O4").Value = "Provincia de Huasco"
$ws.Range("P4").Value = 1240

$ws.Range("D5").Value = 44496
$ws.Range("H5").Value = "Perfection"
$ws.Range("J5").Value = 250
$ws.Range("K5").Value = 14000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 14520
$ws.Range("N5").Value = "$/malla 25 kilos"
$ws.Range("O5").Value = "Provincia de Huasco"
$ws.Range("P5").Value = 581

$ws.Range("D6").Value = 44454
$ws.Range("H6").Value = "Perfection"
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 36000
$ws.Range("L6").Value = 38000
$ws.Range("M6").Value = 37000
$ws.Range("N6").Value = "$/malla 25 kilos"
$ws.Range("O6").Value = "Provincia de Limarí"
$ws.Range("P6").Value = 1480

$ws.Range("D7").Value = 44399
$ws.Range("H7").Value = "Perfection"
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 39000
$ws.Range("L7").Value = 40000
$ws.Range("M7").Value = 39600
$ws.Range("N7").Value = "$/malla 25 kilos"
$ws.Range("O7").Value = "Provincia de Huasco"
$ws.Range("P7").Value = 1584

$ws.Range("D8").Value = 44482
$ws.Range("H8").Value = "Perfection"
$ws.Range("J8").Value = 130
$ws.Range("K8").Value = 24000
$ws.Range("L8").Value = 25000
$ws.Range("M8").Value = 24385
$ws.Range("N8").Value = "$/saco 25 kilos"
$ws.Range("O8").Value = "Región de O'Higgins"
$ws.Range("P8").Value = 975

$ws.Range("D9").Value = 44342
$ws.Range("H9").Value = "Perfection"
$ws.Range("J9").Value = 60
$ws.Range("K9").Value = 30000
$ws.Range("L9").Value = 32000
$ws.Range("M9").Value = 31000
$ws.Range("N9").Value = "$/malla 25 kilos"
$ws.Range("O9").Value = "Provincia de Limarí"
$ws.Range("P9").Value = 1240

